$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.052105665206909
$ws.Range("B1").Value = 3.364169836044312
$ws.Range("C1").Value = 1.991514563560486
$ws.Range("D1").Value = 1.520692944526672
$ws.Range("E1").Value = 1.359121441841125
